# "Added mathhelp, improved output printing"
#
# On the "Assignment 4" sheet, two raw inputs were edited:
#   - L2 (NMC / "a" mark)   50 -> 0
#   - G3 (Language Models / "b" mark) 100 -> 33
# Every other changed cell (B2, D2, N2, B3, D3, I3, D5, N5, I6) is a
# formula whose cached value simply ripples from those two inputs, so
# Excel's automatic recalculation reproduces them once the inputs change.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Assignment 4")
$ws.Activate()

$ws.Range("L2").Value = 0
$ws.Range("G3").Value = 33

# Leave the selection where the user last clicked after making the edits.
$ws.Range("L4").Select()
